$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization "đại" -> "Đại" and "TUQ." -> "T."
$ws.Range("A7").Value = "Trường Đại học FPT xác nhận"
$ws.Range("B27").Value = "Tại trường Đại học FPT."
$ws.Range("F31").Value = "     T. VIỆN TRƯỞNG"

# Update view state (top-left cell + selection)
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("I29").Select() | Out-Null
